# Fixed update to excel issue
# 1. Rename header in "Weekly Quantity" sheet (B1): "Requested quantity" -> "Weekly_PO_Qty"
# 2. Rename header in "Monthly Trend" sheet (B1): "Requested quantity" -> "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" sheet with ds / PO_Forecast / yhat_lower / yhat_upper columns

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Copy header formatting (bold, centered, bordered) from the Weekly Quantity header row
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-number-format from an existing date cell in column A
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A46").PasteSpecial(-4122)

# Header values
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Data rows

$newSheet.Range("A2").Value = 44948.99999999999
$newSheet.Range("B2").Value = 124
$newSheet.Range("C2").Value = -76.84605974825357
$newSheet.Range("D2").Value = 338.5113811113455
$newSheet.Range("A3").Value = 44955.99999999999
$newSheet.Range("B3").Value = 123
$newSheet.Range("C3").Value = -86.26067854119609
$newSheet.Range("D3").Value = 326.9242550845429
$newSheet.Range("A4").Value = 44962.99999999999
$newSheet.Range("B4").Value = 123
$newSheet.Range("C4").Value = -96.74459733588594
$newSheet.Range("D4").Value = 335.5913599603011
$newSheet.Range("A5").Value = 44969.99999999999
$newSheet.Range("B5").Value = 122
$newSheet.Range("C5").Value = -76.68405490100226
$newSheet.Range("D5").Value = 320.1033947606118
$newSheet.Range("A6").Value = 44976.99999999999
$newSheet.Range("B6").Value = 122
$newSheet.Range("C6").Value = -91.82916880286858
$newSheet.Range("D6").Value = 330.9439037766107
$newSheet.Range("A7").Value = 44983.99999999999
$newSheet.Range("B7").Value = 121
$newSheet.Range("C7").Value = -89.55202210515212
$newSheet.Range("D7").Value = 330.3772054680571
$newSheet.Range("A8").Value = 44990.99999999999
$newSheet.Range("B8").Value = 120
$newSheet.Range("C8").Value = -110.1584177063978
$newSheet.Range("D8").Value = 317.7606563735312
$newSheet.Range("A9").Value = 44997.99999999999
$newSheet.Range("B9").Value = 120
$newSheet.Range("C9").Value = -93.55372250633646
$newSheet.Range("D9").Value = 319.6237799651229
$newSheet.Range("A10").Value = 45004.99999999999
$newSheet.Range("B10").Value = 119
$newSheet.Range("C10").Value = -88.444351284481
$newSheet.Range("D10").Value = 328.9034414929499
$newSheet.Range("A11").Value = 45018.99999999999
$newSheet.Range("B11").Value = 118
$newSheet.Range("C11").Value = -85.87369324569968
$newSheet.Range("D11").Value = 316.686767044467
$newSheet.Range("A12").Value = 45025.99999999999
$newSheet.Range("B12").Value = 118
$newSheet.Range("C12").Value = -87.29957445314554
$newSheet.Range("D12").Value = 341.2162548556786
$newSheet.Range("A13").Value = 45032.99999999999
$newSheet.Range("B13").Value = 117
$newSheet.Range("C13").Value = -100.642278505436
$newSheet.Range("D13").Value = 324.0379495028477
$newSheet.Range("A14").Value = 45039.99999999999
$newSheet.Range("B14").Value = 117
$newSheet.Range("C14").Value = -93.43351520046348
$newSheet.Range("D14").Value = 319.2509886170785
$newSheet.Range("A15").Value = 45046.99999999999
$newSheet.Range("B15").Value = 116
$newSheet.Range("C15").Value = -102.1032088197203
$newSheet.Range("D15").Value = 332.9292681232612
$newSheet.Range("A16").Value = 45053.99999999999
$newSheet.Range("B16").Value = 115
$newSheet.Range("C16").Value = -104.0453531549205
$newSheet.Range("D16").Value = 316.0568592692078
$newSheet.Range("A17").Value = 45060.99999999999
$newSheet.Range("B17").Value = 115
$newSheet.Range("C17").Value = -94.46490018719729
$newSheet.Range("D17").Value = 328.1664886406542
$newSheet.Range("A18").Value = 45067.99999999999
$newSheet.Range("B18").Value = 114
$newSheet.Range("C18").Value = -105.4350561320652
$newSheet.Range("D18").Value = 318.5249299906797
$newSheet.Range("A19").Value = 45074.99999999999
$newSheet.Range("B19").Value = 114
$newSheet.Range("C19").Value = -96.80060939583814
$newSheet.Range("D19").Value = 320.132270271609
$newSheet.Range("A20").Value = 45081.99999999999
$newSheet.Range("B20").Value = 113
$newSheet.Range("C20").Value = -101.1348390131342
$newSheet.Range("D20").Value = 324.4603380239112
$newSheet.Range("A21").Value = 45088.99999999999
$newSheet.Range("B21").Value = 113
$newSheet.Range("C21").Value = -92.76785186171725
$newSheet.Range("D21").Value = 334.0135688530568
$newSheet.Range("A22").Value = 45116.99999999999
$newSheet.Range("B22").Value = 110
$newSheet.Range("C22").Value = -98.64840998862535
$newSheet.Range("D22").Value = 321.313633156797
$newSheet.Range("A23").Value = 45130.99999999999
$newSheet.Range("B23").Value = 109
$newSheet.Range("C23").Value = -108.1983232428198
$newSheet.Range("D23").Value = 323.7810885432965
$newSheet.Range("A24").Value = 45137.99999999999
$newSheet.Range("B24").Value = 109
$newSheet.Range("C24").Value = -96.04739354349871
$newSheet.Range("D24").Value = 321.0163686567577
$newSheet.Range("A25").Value = 45144.99999999999
$newSheet.Range("B25").Value = 108
$newSheet.Range("C25").Value = -97.2764900802232
$newSheet.Range("D25").Value = 324.2357379466874
$newSheet.Range("A26").Value = 45151.99999999999
$newSheet.Range("B26").Value = 108
$newSheet.Range("C26").Value = -116.394832997366
$newSheet.Range("D26").Value = 316.9659718290417
$newSheet.Range("A27").Value = 45158.99999999999
$newSheet.Range("B27").Value = 107
$newSheet.Range("C27").Value = -100.9995617612699
$newSheet.Range("D27").Value = 320.4987491109408
$newSheet.Range("A28").Value = 45165.99999999999
$newSheet.Range("B28").Value = 106
$newSheet.Range("C28").Value = -96.52167984752026
$newSheet.Range("D28").Value = 301.5058422564326
$newSheet.Range("A29").Value = 45172.99999999999
$newSheet.Range("B29").Value = 106
$newSheet.Range("C29").Value = -93.11258382351402
$newSheet.Range("D29").Value = 313.4449462713008
$newSheet.Range("A30").Value = 45179.99999999999
$newSheet.Range("B30").Value = 105
$newSheet.Range("C30").Value = -112.556996605349
$newSheet.Range("D30").Value = 298.0075452183194
$newSheet.Range("A31").Value = 45200.99999999999
$newSheet.Range("B31").Value = 104
$newSheet.Range("C31").Value = -102.9066431790916
$newSheet.Range("D31").Value = 306.6654447777082
$newSheet.Range("A32").Value = 45207.99999999999
$newSheet.Range("B32").Value = 103
$newSheet.Range("C32").Value = -88.69269291686815
$newSheet.Range("D32").Value = 306.9304617569896
$newSheet.Range("A33").Value = 45214.99999999999
$newSheet.Range("B33").Value = 103
$newSheet.Range("C33").Value = -103.4041313137168
$newSheet.Range("D33").Value = 308.8037991694947
$newSheet.Range("A34").Value = 45221.99999999999
$newSheet.Range("B34").Value = 102
$newSheet.Range("C34").Value = -113.1419391679785
$newSheet.Range("D34").Value = 309.1181720909646
$newSheet.Range("A35").Value = 45235.99999999999
$newSheet.Range("B35").Value = 101
$newSheet.Range("C35").Value = -100.277018746593
$newSheet.Range("D35").Value = 318.8248146634285
$newSheet.Range("A36").Value = 45249.99999999999
$newSheet.Range("B36").Value = 100
$newSheet.Range("C36").Value = -100.8774133289383
$newSheet.Range("D36").Value = 319.334753934934
$newSheet.Range("A37").Value = 45256.99999999999
$newSheet.Range("B37").Value = 99
$newSheet.Range("C37").Value = -107.4747673895341
$newSheet.Range("D37").Value = 313.2166496105365
$newSheet.Range("A38").Value = 45312.99999999999
$newSheet.Range("B38").Value = 95
$newSheet.Range("C38").Value = -118.3565632672228
$newSheet.Range("D38").Value = 307.0658434744197
$newSheet.Range("A39").Value = 45319.99999999999
$newSheet.Range("B39").Value = 94
$newSheet.Range("C39").Value = -121.1355715754816
$newSheet.Range("D39").Value = 305.595222133179
$newSheet.Range("A40").Value = 45326.99999999999
$newSheet.Range("B40").Value = 94
$newSheet.Range("C40").Value = -115.8910952302836
$newSheet.Range("D40").Value = 296.2838216344215
$newSheet.Range("A41").Value = 45333.99999999999
$newSheet.Range("B41").Value = 93
$newSheet.Range("C41").Value = -125.1747146487971
$newSheet.Range("D41").Value = 311.5305223157795
$newSheet.Range("A42").Value = 45340.99999999999
$newSheet.Range("B42").Value = 92
$newSheet.Range("C42").Value = -102.9467813166758
$newSheet.Range("D42").Value = 308.9563655129838
$newSheet.Range("A43").Value = 45347.99999999999
$newSheet.Range("B43").Value = 92
$newSheet.Range("C43").Value = -125.7362351038202
$newSheet.Range("D43").Value = 293.619297796478
$newSheet.Range("A44").Value = 45354.99999999999
$newSheet.Range("B44").Value = 91
$newSheet.Range("C44").Value = -129.6634217478771
$newSheet.Range("D44").Value = 306.9476122548942
$newSheet.Range("A45").Value = 45361.99999999999
$newSheet.Range("B45").Value = 91
$newSheet.Range("C45").Value = -121.9629643715211
$newSheet.Range("D45").Value = 309.7841457348561
$newSheet.Range("A46").Value = 45368.99999999999
$newSheet.Range("B46").Value = 90
$newSheet.Range("C46").Value = -125.6232930469978
$newSheet.Range("D46").Value = 298.6174842711009
